{"js": "// Update the date line and the twenty-five \"three-digit x one-digit\" practice\n// answers in the table to the new day's values.\nconst replacements = [\n  [\"2024-12-16 Monday\", \"2024-12-17 Tuesday\"],\n  [\"444\u00d74=1776\", \"376\u00d73=1128\"],\n  [\"994\u00d75=4970\", \"429\u00d72=858\"],\n  [\"713\u00d72=1426\", \"840\u00d77=5880\"],\n  [\"788\u00d75=3940\", \"907\u00d74=3628\"],\n  [\"779\u00d75=3895\", \"246\u00d76=1476\"],\n  [\"395\u00d79=3555\", \"603\u00d75=3015\"],\n  [\"896\u00d79=8064\", \"466\u00d78=3728\"],\n  [\"677\u00d79=6093\", \"559\u00d75=2795\"],\n  [\"843\u00d74=3372\", \"666\u00d72=1332\"],\n  [\"755\u00d76=4530\", \"132\u00d73=396\"],\n  [\"456\u00d77=3192\", \"978\u00d76=5868\"],\n  [\"937\u00d77=6559\", \"845\u00d77=5915\"],\n  [\"587\u00d77=4109\", \"214\u00d75=1070\"],\n  [\"153\u00d77=1071\", \"271\u00d78=2168\"],\n  [\"428\u00d75=2140\", \"855\u00d77=5985\"],\n  [\"664\u00d72=1328\", \"989\u00d76=5934\"],\n  [\"956\u00d77=6692\", \"461\u00d77=3227\"],\n  [\"997\u00d72=1994\", \"474\u00d77=3318\"],\n  [\"910\u00d77=6370\", \"959\u00d78=7672\"],\n  [\"900\u00d77=6300\", \"389\u00d79=3501\"],\n  [\"854\u00d79=7686\", \"822\u00d78=6576\"],\n  [\"471\u00d78=3768\", \"530\u00d73=1590\"],\n  [\"350\u00d75=1750\", \"840\u00d76=5040\"],\n  [\"759\u00d74=3036\", \"223\u00d73=669\"],\n  [\"489\u00d75=2445\", \"941\u00d77=6587\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the twenty-five \"three-digit x one-digit\" practice\n# answers in the table to the new day's values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2024-12-16 Monday\", \"2024-12-17 Tuesday\"),\n  @(\"444\u00d74=1776\", \"376\u00d73=1128\"),\n  @(\"994\u00d75=4970\", \"429\u00d72=858\"),\n  @(\"713\u00d72=1426\", \"840\u00d77=5880\"),\n  @(\"788\u00d75=3940\", \"907\u00d74=3628\"),\n  @(\"779\u00d75=3895\", \"246\u00d76=1476\"),\n  @(\"395\u00d79=3555\", \"603\u00d75=3015\"),\n  @(\"896\u00d79=8064\", \"466\u00d78=3728\"),\n  @(\"677\u00d79=6093\", \"559\u00d75=2795\"),\n  @(\"843\u00d74=3372\", \"666\u00d72=1332\"),\n  @(\"755\u00d76=4530\", \"132\u00d73=396\"),\n  @(\"456\u00d77=3192\", \"978\u00d76=5868\"),\n  @(\"937\u00d77=6559\", \"845\u00d77=5915\"),\n  @(\"587\u00d77=4109\", \"214\u00d75=1070\"),\n  @(\"153\u00d77=1071\", \"271\u00d78=2168\"),\n  @(\"428\u00d75=2140\", \"855\u00d77=5985\"),\n  @(\"664\u00d72=1328\", \"989\u00d76=5934\"),\n  @(\"956\u00d77=6692\", \"461\u00d77=3227\"),\n  @(\"997\u00d72=1994\", \"474\u00d77=3318\"),\n  @(\"910\u00d77=6370\", \"959\u00d78=7672\"),\n  @(\"900\u00d77=6300\", \"389\u00d79=3501\"),\n  @(\"854\u00d79=7686\", \"822\u00d78=6576\"),\n  @(\"471\u00d78=3768\", \"530\u00d73=1590\"),\n  @(\"350\u00d75=1750\", \"840\u00d76=5040\"),\n  @(\"759\u00d74=3036\", \"223\u00d73=669\"),\n  @(\"489\u00d75=2445\", \"941\u00d77=6587\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $rng = $d.Content\n  $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    Write-Output \"NOT FOUND: $oldText\"\n  }\n}\n"}
